# Adds a new "PersonID -> Name" lookup table in columns F:G (header in row 5,
# data in rows 17-22), reusing the box-border look already used by the
# Person1..Person6 / TaskA..TaskC table in columns A:D.
#
# -4122 == xlPasteFormats (PasteSpecial paste-what argument): copying an
# existing cell and pasting only its formats lets us reuse the exact same
# border/font combination (and, for two new combinations, build on top of an
# existing one) so the resulting style indices match what Excel itself would
# have produced/deduplicated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- "middle" rows of the new table (F18:F21 / G18:G21) ---
# Same style as the existing A8:A15 (bold, left border) / D8:D15 (right
# border) cells of the Person list / task table.
$ws.Range("A7").Copy()
$ws.Range("F18:F21").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("G18:G21").PasteSpecial(-4122)

$ws.Range("F18").Value = "Person2"
$ws.Range("G18").Value = "Mark"
$ws.Range("F19").Value = "Person3"
$ws.Range("G19").Value = "Zihan"
$ws.Range("F20").Value = "Person4"
$ws.Range("G20").Value = "Yifan"
$ws.Range("F21").Value = "Person5"
$ws.Range("G21").Value = "Saad"

# --- bottom row of the new table (F22 / G22) ---
# Same style as A16 (bold, left+bottom border) / D16 (right+bottom border).
$ws.Range("A16").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("G22").PasteSpecial(-4122)

$ws.Range("F22").Value = "Person6"
$ws.Range("G22").Value = "Nicky"

# --- top row of the new table (F17 / G17) ---
# F17 needs a NEW style: bold + left+top border. Start from A4's left+top
# border (not bold) and add Bold on top -> engine creates/reuses the new
# combined style.
# G17 reuses the existing right+top border style already used by D4.
$ws.Range("A4").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F17").Font.Bold = $true
$ws.Range("D4").Copy()
$ws.Range("G17").PasteSpecial(-4122)

$ws.Range("F17").Value = "Person1"
$ws.Range("G17").Value = "Arun"

# --- header row (F5 / G5) ---
# F5 reuses the bold+left-border style (same as F18 etc. above).
# G5 needs a NEW style: bold + right border + right-aligned text. Start from
# D7's right-border style (not bold) and add Bold + right alignment.
$ws.Range("A7").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Font.Bold = $true
$ws.Range("G5").HorizontalAlignment = -4152

$ws.Range("F5").Value = "Week"
$ws.Range("G5").Value = "Person"

# Match the author's final cursor position.
$ws.Range("H7").Select()
